# "Updated_Format_1 (Text Wrapper+Coloana B mai mica)"
#
# 1. Prefix the "Numele și Prenumele:" label (A3) with two leading spaces.
# 2. Recolor the grey accent font (style index 2, used by A1 / the word-count
#    cells) from FF808080 to FFE3E3E3, and center those cells.
# 3. Turn on text-wrapping for the headword/definition blocks (style 6,
#    additionally centered horizontally) and the alternating answer-row
#    fills (styles 7, 8, 9).
# 4. Narrow column B from 250.71 to 200.71 characters.
# 5. Stamp a literal "10" (same look as the other word-count cells) into the
#    A column of the 3rd answer row of every block (A7, A19, A31, ... A115).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Label text tweak -------------------------------------------------
$ws.Range("A3").Value = "  Numele și Prenumele:"

# --- 2. Grey "word count" cells: recolor + center -----------------------
# Set up one reference cell fully, then fan the resulting format out to its
# siblings via copy/paste-special so we don't leave orphaned intermediate
# styles behind in styles.xml.
$countCells = @("A1","A6","A18","A30","A42","A54","A66","A78","A90","A102","A114")
$refCount = $ws.Range($countCells[0])
$refCount.Font.Color = 14935011          # RGB(227,227,227) = FFE3E3E3
$refCount.HorizontalAlignment = -4108    # xlCenter
$refCount.VerticalAlignment = -4108      # xlCenter
$refCount.Copy()
foreach ($c in $countCells) {
    if ($c -ne $countCells[0]) {
        $ws.Range($c).PasteSpecial(-4122)   # xlPasteFormats
    }
}

# --- 3. Wrap text on the colored blocks ----------------------------------
$headwordCells = @("A5","A17","A29","A41","A53","A65","A77","A89","A101","A113")
$refHead = $ws.Range($headwordCells[0])
$refHead.HorizontalAlignment = -4108    # xlCenter
$refHead.VerticalAlignment = -4108      # xlCenter
$refHead.WrapText = $true
$refHead.Copy()
foreach ($c in $headwordCells) {
    if ($c -ne $headwordCells[0]) {
        $ws.Range($c).PasteSpecial(-4122)
    }
}

$wrapGroups = @(
    @("B5","B17","B29","B41","B53","B65","B77","B89","B101","B113"),
    @("B6","B18","B30","B42","B54","B66","B78","B90","B102","B114","B116","B118","B120","B122"),
    @("B7","B19","B31","B43","B55","B67","B79","B91","B103","B115","B117","B119","B121","B123")
)
foreach ($grp in $wrapGroups) {
    $refCell = $ws.Range($grp[0])
    $refCell.WrapText = $true
    $refCell.Copy()
    foreach ($c in $grp) {
        if ($c -ne $grp[0]) {
            $ws.Range($c).PasteSpecial(-4122)
        }
    }
}

# --- 4. Narrow column B ---------------------------------------------------
# Target author width is 200.7109375 chars; the host quantizes ColumnWidth to
# whole pixels, so feed it the input that lands on the nearest achievable
# grid point (200.6667) instead of the input value itself (which rounds the
# other way, to 201.5).
$ws.Columns("B").ColumnWidth = 199.8

# --- 5. Add the "10" marker to the 3rd answer row of every block --------
$markerRows = @(7,19,31,43,55,67,79,91,103,115)
$ws.Range("A1").Copy()
foreach ($r in $markerRows) {
    $target = $ws.Cells.Item($r, 1)
    $target.Value = 10
    $target.PasteSpecial(-4122)   # xlPasteFormats
}
